$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-valued ranges that must stay as literal text (dates, numeric-looking
# codes with leading zeros, long document numbers, decimal amounts stored
# as text in the source data). Force text format so COM doesn't silently
# reinterpret them as dates/numbers.
$textRanges = @(
  "B2","B3","B4","B5",
  "E2","E3","E4","E5",
  "F2","F3","F4","F5",
  "G2","G3","G4","G5",
  "I2","I5",
  "J2","J3","J4","J5",
  "K2","K5",
  "M2","M3","M4","M5",
  "N2","N3","N4","N5",
  "O2","O3","O4","O5",
  "Q2","Q3","Q4","Q5",
  "R2",
  "S2","S3","S4","S5"
)
foreach ($addr in $textRanges) {
  $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("B2").Value = "2024-07-16"
$ws.Range("D2").Value = 88206
$ws.Range("E2").Value = "100094507"
$ws.Range("F2").Value = "BR0026113"
$ws.Range("G2").Value = "MANGUINHOS ADM DE BENS E CONSULTORI"
$ws.Range("I2").Value = "Remessa"
$ws.Range("J2").Value = "41240776881093000172550020000882061000945072"
$ws.Range("K2").Value = "BR0105"
$ws.Range("M2").Value = "48699418000165"
$ws.Range("N2").Value = "04557"
$ws.Range("O2").Value = "RJ"
$ws.Range("Q2").Value = "280.26"
$ws.Range("R2").Value = "70.07"
$ws.Range("S2").Value = "350.33"

# Row 3
$ws.Range("B3").Value = "2024-07-16"
$ws.Range("D3").Value = 88208
$ws.Range("E3").Value = "100094535"
$ws.Range("F3").Value = "BR0010977"
$ws.Range("G3").Value = "AMB EMPREENDIMENTOS IMOBIL. LTDA"
$ws.Range("J3").Value = "41240776881093000172550020000882081000945352"
$ws.Range("M3").Value = "26354329000144"
$ws.Range("N3").Value = "02404"
$ws.Range("O3").Value = "SC"
$ws.Range("Q3").Value = "2125.02"
$ws.Range("S3").Value = "2125.02"

# Row 4
$ws.Range("B4").Value = "2024-07-16"
$ws.Range("D4").Value = 88211
$ws.Range("E4").Value = "100094538"
$ws.Range("F4").Value = "BR0025869"
$ws.Range("G4").Value = "CONDOMINIO SOBERANE RESIDENCE, CORP"
$ws.Range("J4").Value = "41240776881093000172550020000882111000945382"
$ws.Range("M4").Value = "32581733000153"
$ws.Range("N4").Value = "02603"
$ws.Range("O4").Value = "AM"
$ws.Range("Q4").Value = "793.61"
$ws.Range("S4").Value = "793.61"

# Row 5
$ws.Range("B5").Value = "2024-07-16"
$ws.Range("D5").Value = 88212
$ws.Range("E5").Value = "100094541"
$ws.Range("F5").Value = "BR0015419"
$ws.Range("G5").Value = "SPE SAUDE PRIMARIA BH S/A"
$ws.Range("I5").Value = "Venda com pedido"
$ws.Range("J5").Value = "41240776881093000172550020000882121000945410"
$ws.Range("K5").Value = "BR0101"
$ws.Range("M5").Value = "23921007000141"
$ws.Range("N5").Value = "06200"
$ws.Range("O5").Value = "MG"
$ws.Range("Q5").Value = "4397.29"
$ws.Range("S5").Value = "4397.29"
